$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9930.458000000001
$ws.Range("J116").Value = 3523.625
$ws.Range("L116").Value = 3523.625
$ws.Range("N116").Value = -10407.625
$ws.Range("H125").Value = 2861.6155
$ws.Range("I125").Value = 1666
$ws.Range("J125").Value = 3886.4285
$ws.Range("K125").Value = 14994
$ws.Range("L125").Value = 34977.8565
$ws.Range("M125").Value = -12534
$ws.Range("N125").Value = -39897.8565
$ws.Range("H137").Value = 4627.4287
$ws.Range("I137").Value = 11201
$ws.Range("J137").Value = 1998
$ws.Range("K137").Value = 33603
$ws.Range("L137").Value = 5994
$ws.Range("M137").Value = -31053
$ws.Range("N137").Value = -11094
$ws.Range("H138").Value = 3287.9082
$ws.Range("I138").Value = 1780.7
$ws.Range("J138").Value = 3952.853
$ws.Range("K138").Value = 5342.1
$ws.Range("L138").Value = 11858.559
$ws.Range("M138").Value = -202.1000000000004
$ws.Range("N138").Value = -22138.559
$ws.Range("H140").Value = 199965
$ws.Range("J140").Value = 199965
$ws.Range("L140").Value = 199965
$ws.Range("N140").Value = -210325
$ws.Range("H141").Value = 6789.727
$ws.Range("I141").Value = 6214.6665
$ws.Range("K141").Value = 18643.9995
$ws.Range("M141").Value = -13463.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18182738
$ws.Range("I32").Value = 19231454
$ws.Range("K32").Value = 19231454
$ws.Range("M32").Value = -19231167
$ws.Range("H45").Value = 3935.8
$ws.Range("J45").Value = 4583.1665
$ws.Range("L45").Value = 4583.1665
$ws.Range("N45").Value = -5337.1665
$ws.Range("H61").Value = 2211.5625
$ws.Range("I61").Value = 2006.1724
$ws.Range("J61").Value = 4197
$ws.Range("K61").Value = 2006.1724
$ws.Range("L61").Value = 4197
$ws.Range("M61").Value = -1794.1724
$ws.Range("N61").Value = -4621
$ws.Range("H74").Value = 1926.1666
$ws.Range("I74").Value = 2516.5454
$ws.Range("J74").Value = 998.4286
$ws.Range("K74").Value = 2516.5454
$ws.Range("L74").Value = 998.4286
$ws.Range("M74").Value = -1642.5454
$ws.Range("N74").Value = -2746.4286
$ws.Range("H77").Value = 1926.1666
$ws.Range("I77").Value = 2516.5454
$ws.Range("J77").Value = 998.4286
$ws.Range("K77").Value = 12582.727
$ws.Range("L77").Value = 4992.143
$ws.Range("M77").Value = -8214.726999999999
$ws.Range("N77").Value = -13728.143
$ws.Range("H132").Value = 2042.6
$ws.Range("I132").Value = 1966.7693
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5900.3079
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -3370.3079
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 2211.5625
$ws.Range("I136").Value = 2006.1724
$ws.Range("J136").Value = 4197
$ws.Range("K136").Value = 6018.5172
$ws.Range("L136").Value = 12591
$ws.Range("M136").Value = -3468.5172
$ws.Range("N136").Value = -17691

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24619.363
$ws.Range("I82").Value = 20581.3
$ws.Range("J82").Value = 65000
$ws.Range("K82").Value = 20581.3
$ws.Range("L82").Value = 65000
$ws.Range("M82").Value = -20198.3
$ws.Range("N82").Value = -65766
$ws.Range("H85").Value = 24619.363
$ws.Range("I85").Value = 20581.3
$ws.Range("J85").Value = 65000
$ws.Range("K85").Value = 20581.3
$ws.Range("L85").Value = 65000
$ws.Range("M85").Value = -19255.3
$ws.Range("N85").Value = -67652
$ws.Range("H86").Value = 2871.75
$ws.Range("I86").Value = 2646
$ws.Range("K86").Value = 2646
$ws.Range("M86").Value = -1523
$ws.Range("H89").Value = 2871.75
$ws.Range("I89").Value = 2646
$ws.Range("K89").Value = 13230
$ws.Range("M89").Value = -7614
$ws.Range("H94").Value = 536.2273
$ws.Range("I94").Value = 557
$ws.Range("J94").Value = 100
$ws.Range("K94").Value = 557
$ws.Range("L94").Value = 100
$ws.Range("M94").Value = -106
$ws.Range("N94").Value = -1002
$ws.Range("H99").Value = 3521.2307
$ws.Range("I99").Value = 3635.6
$ws.Range("J99").Value = 3449.75
$ws.Range("K99").Value = 3635.6
$ws.Range("L99").Value = 3449.75
$ws.Range("M99").Value = -2137.6
$ws.Range("N99").Value = -6445.75
$ws.Range("H134").Value = 1698.9048
$ws.Range("I134").Value = 1226.8889
$ws.Range("K134").Value = 3680.6667
$ws.Range("M134").Value = -1145.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1713.2424
$ws.Range("I31").Value = 857
$ws.Range("J31").Value = 2569.4849
$ws.Range("K31").Value = 857
$ws.Range("L31").Value = 2569.4849
$ws.Range("M31").Value = -562
$ws.Range("N31").Value = -3159.4849
$ws.Range("H34").Value = 1713.2424
$ws.Range("I34").Value = 857
$ws.Range("J34").Value = 2569.4849
$ws.Range("K34").Value = 857
$ws.Range("L34").Value = 2569.4849
$ws.Range("M34").Value = -655
$ws.Range("N34").Value = -2973.4849
$ws.Range("H58").Value = 2832.175
$ws.Range("I58").Value = 1497
$ws.Range("K58").Value = 1497
$ws.Range("M58").Value = -1294
$ws.Range("H104").Value = 79975
$ws.Range("J104").Value = 79975
$ws.Range("L104").Value = 79975
$ws.Range("N104").Value = -85217
$ws.Range("H132").Value = 5767.5884
$ws.Range("I132").Value = 5753.25
$ws.Range("K132").Value = 17259.75
$ws.Range("M132").Value = -14729.75
$ws.Range("H134").Value = 4985.1113
$ws.Range("I134").Value = 5783.933
$ws.Range("J134").Value = 991
$ws.Range("K134").Value = 17351.799
$ws.Range("L134").Value = 2973
$ws.Range("M134").Value = -14816.799
$ws.Range("N134").Value = -8043
$ws.Range("H136").Value = 2832.175
$ws.Range("I136").Value = 1497
$ws.Range("K136").Value = 4491
$ws.Range("M136").Value = -1941
$ws.Range("H141").Value = 652835.7
$ws.Range("J141").Value = 652835.7
$ws.Range("L141").Value = 652835.7
$ws.Range("N141").Value = -663195.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 173.33333
$ws.Range("J16").Value = 172
$ws.Range("L16").Value = 516
$ws.Range("N16").Value = -862
$ws.Range("H56").Value = 5149.0415
$ws.Range("I56").Value = 5149.0415
$ws.Range("K56").Value = 5149.0415
$ws.Range("M56").Value = -4619.0415
$ws.Range("H134").Value = 3113.6667
$ws.Range("I134").Value = 3113.6667
$ws.Range("K134").Value = 9341.000100000001
$ws.Range("M134").Value = -4271.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 20833.334
$ws.Range("I40").Value = 17500
$ws.Range("J40").Value = 27500
$ws.Range("K40").Value = 17500
$ws.Range("L40").Value = 27500
$ws.Range("M40").Value = -17349
$ws.Range("N40").Value = -27802
$ws.Range("H132").Value = 2404.875
$ws.Range("I132").Value = 1765.4445
$ws.Range("K132").Value = 5296.333500000001
$ws.Range("M132").Value = -2766.333500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4497.25
$ws.Range("I16").Value = 4329.6665
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 4329.6665
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -4159.6665
$ws.Range("N16").Value = -5340
$ws.Range("H46").Value = 6666.567
$ws.Range("J46").Value = 7223.2964
$ws.Range("L46").Value = 7223.2964
$ws.Range("N46").Value = -7599.2964
$ws.Range("H136").Value = 53135.285
$ws.Range("I136").Value = 61324.668
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 183974.004
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -181424.004
$ws.Range("N136").Value = -17097
$ws.Range("H138").Value = 172241.5
$ws.Range("J138").Value = 172241.5
$ws.Range("L138").Value = 172241.5
$ws.Range("N138").Value = -182521.5
$ws.Range("H140").Value = 89749.75
$ws.Range("J140").Value = 86333
$ws.Range("L140").Value = 86333
$ws.Range("N140").Value = -96693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 381.22726
$ws.Range("I107").Value = 389.38095
$ws.Range("K107").Value = 1168.14285
$ws.Range("M107").Value = 751.85715
$ws.Range("H132").Value = 2470.2222
$ws.Range("I132").Value = 2468.4707
$ws.Range("K132").Value = 7405.4121
$ws.Range("M132").Value = -4875.4121
$ws.Range("H136").Value = 59645.777
$ws.Range("I136").Value = 3519.0908
$ws.Range("J136").Value = 147844.86
$ws.Range("K136").Value = 10557.2724
$ws.Range("L136").Value = 443534.58
$ws.Range("M136").Value = -8007.2724
$ws.Range("N136").Value = -448634.58
